# Update "想去人数" (people interested) counts that changed between data pulls.
# Sheet "展览" (index 1) and sheet "全部类型" (index 4) both contain the same
# events, so each updated number needs to be applied in both places.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet updates
$wsExhibition.Range("F2").Value = 5286
$wsExhibition.Range("F4").Value = 10760
$wsExhibition.Range("F6").Value = 574
$wsExhibition.Range("F8").Value = 187
$wsExhibition.Range("F9").Value = 900

# 全部类型 sheet updates (same events, different row positions)
$wsAllTypes.Range("F4").Value = 5286
$wsAllTypes.Range("F7").Value = 10760
$wsAllTypes.Range("F9").Value = 574
$wsAllTypes.Range("F13").Value = 187
$wsAllTypes.Range("F14").Value = 900
